# Auto-generated script to update Unicorn_Profits market-data values
# per commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4371.815
$ws.Range("I64").Value = 3927.5557
$ws.Range("J64").Value = 4593.9443
$ws.Range("K64").Value = 3927.5557
$ws.Range("L64").Value = 4593.9443
$ws.Range("M64").Value = -3679.5557
$ws.Range("N64").Value = -5089.9443
$ws.Range("H67").Value = 4371.815
$ws.Range("I67").Value = 3927.5557
$ws.Range("J67").Value = 4593.9443
$ws.Range("K67").Value = 3927.5557
$ws.Range("L67").Value = 4593.9443
$ws.Range("M67").Value = -3069.5557
$ws.Range("N67").Value = -6309.9443
$ws.Range("H137").Value = 73338.53
$ws.Range("I137").Value = 107808.1
$ws.Range("J137").Value = 4399.4
$ws.Range("K137").Value = 323424.3
$ws.Range("L137").Value = 13198.2
$ws.Range("M137").Value = -320874.3
$ws.Range("N137").Value = -18298.2

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 40361.645
$ws.Range("I134").Value = 68501.734
$ws.Range("J134").Value = 7892.3076
$ws.Range("K134").Value = 205505.202
$ws.Range("L134").Value = 23676.9228
$ws.Range("M134").Value = -202970.202
$ws.Range("N134").Value = -28746.9228

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181
$ws.Range("H31").Value = 3410.1025
$ws.Range("I31").Value = 2203.8096
$ws.Range("J31").Value = 4817.4443
$ws.Range("K31").Value = 2203.8096
$ws.Range("L31").Value = 4817.4443
$ws.Range("M31").Value = -1908.8096
$ws.Range("N31").Value = -5407.4443
$ws.Range("H34").Value = 3410.1025
$ws.Range("I34").Value = 2203.8096
$ws.Range("J34").Value = 4817.4443
$ws.Range("K34").Value = 2203.8096
$ws.Range("L34").Value = 4817.4443
$ws.Range("M34").Value = -2001.8096
$ws.Range("N34").Value = -5221.4443
$ws.Range("H58").Value = 2277.8462
$ws.Range("I58").Value = 1623.8462
$ws.Range("J58").Value = 2931.8462
$ws.Range("K58").Value = 1623.8462
$ws.Range("L58").Value = 2931.8462
$ws.Range("M58").Value = -1420.8462
$ws.Range("N58").Value = -3337.8462
$ws.Range("H62").Value = 2225266.2
$ws.Range("J62").Value = 3307.2307
$ws.Range("L62").Value = 3307.2307
$ws.Range("N62").Value = -4555.2307
$ws.Range("H65").Value = 2225266.2
$ws.Range("J65").Value = 3307.2307
$ws.Range("L65").Value = 16536.1535
$ws.Range("N65").Value = -22776.1535
$ws.Range("H99").Value = 30809.885
$ws.Range("I99").Value = 43426.293
$ws.Range("J99").Value = 3283.182
$ws.Range("K99").Value = 43426.293
$ws.Range("L99").Value = 3283.182
$ws.Range("M99").Value = -41928.293
$ws.Range("N99").Value = -6279.182
$ws.Range("H122").Value = 2538.5454
$ws.Range("I122").Value = 2538.5454
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7615.6362
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5165.6362
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 30809.885
$ws.Range("I126").Value = 43426.293
$ws.Range("J126").Value = 3283.182
$ws.Range("K126").Value = 130278.879
$ws.Range("L126").Value = 9849.545999999998
$ws.Range("M126").Value = -127808.879
$ws.Range("N126").Value = -14789.546
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H136").Value = 2277.8462
$ws.Range("I136").Value = 1623.8462
$ws.Range("J136").Value = 2931.8462
$ws.Range("K136").Value = 4871.5386
$ws.Range("L136").Value = 8795.5386
$ws.Range("M136").Value = -2321.5386
$ws.Range("N136").Value = -13895.5386

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1327.55
$ws.Range("I97").Value = 904.1667
$ws.Range("J97").Value = 1962.625
$ws.Range("K97").Value = 904.1667
$ws.Range("L97").Value = 1962.625
$ws.Range("M97").Value = -408.1667
$ws.Range("N97").Value = -2954.625

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 333.75
$ws.Range("I22").Value = 264
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 264
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = 31
$ws.Range("N22").Value = -1040
$ws.Range("H27").Value = 333.75
$ws.Range("I27").Value = 264
$ws.Range("J27").Value = 450
$ws.Range("K27").Value = 264
$ws.Range("L27").Value = 450
$ws.Range("M27").Value = -157
$ws.Range("N27").Value = -664
$ws.Range("H46").Value = 869.6
$ws.Range("I46").Value = 748.8
$ws.Range("J46").Value = 930
$ws.Range("K46").Value = 748.8
$ws.Range("L46").Value = 930
$ws.Range("M46").Value = -560.8
$ws.Range("N46").Value = -1306
$ws.Range("H55").Value = 280.2857
$ws.Range("I55").Value = 243.14285
$ws.Range("K55").Value = 243.14285
$ws.Range("M55").Value = -70.14285000000001
$ws.Range("H68").Value = 2665.2856
$ws.Range("I68").Value = 1997.6364
$ws.Range("J68").Value = 3399.7
$ws.Range("K68").Value = 1997.6364
$ws.Range("L68").Value = 3399.7
$ws.Range("M68").Value = -1248.6364
$ws.Range("N68").Value = -4897.7
$ws.Range("H71").Value = 2665.2856
$ws.Range("I71").Value = 1997.6364
$ws.Range("J71").Value = 3399.7
$ws.Range("K71").Value = 9988.182000000001
$ws.Range("L71").Value = 16998.5
$ws.Range("M71").Value = -6244.182000000001
$ws.Range("N71").Value = -24486.5
$ws.Range("H93").Value = 1695.0741
$ws.Range("I93").Value = 1708.1666
$ws.Range("J93").Value = 1668.8889
$ws.Range("K93").Value = 1708.1666
$ws.Range("L93").Value = 1668.8889
$ws.Range("M93").Value = -460.1666
$ws.Range("N93").Value = -4164.8889

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4407.8096
$ws.Range("I62").Value = 4314.4
$ws.Range("J62").Value = 4437
$ws.Range("K62").Value = 4314.4
$ws.Range("L62").Value = 4437
$ws.Range("M62").Value = -3690.4
$ws.Range("N62").Value = -5685
$ws.Range("H65").Value = 4407.8096
$ws.Range("I65").Value = 4314.4
$ws.Range("J65").Value = 4437
$ws.Range("K65").Value = 21572
$ws.Range("L65").Value = 22185
$ws.Range("M65").Value = -18452
$ws.Range("N65").Value = -28425
$ws.Range("H132").Value = 38615.18
$ws.Range("I132").Value = 73265.21000000001
$ws.Range("J132").Value = 3965.1428
$ws.Range("K132").Value = 219795.63
$ws.Range("L132").Value = 11895.4284
$ws.Range("M132").Value = -217265.63
$ws.Range("N132").Value = -16955.4284
$ws.Range("H138").Value = 30011.8
$ws.Range("J138").Value = 30011.8
$ws.Range("L138").Value = 30011.8
$ws.Range("N138").Value = -40291.8
